$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This edit rearranges the *content* of a block of observation rows
# (17,18,19,21,22,23,24,25,26,29,30) on sheet "Artfynd": every cell's data
# moves to a different row number per a fixed permutation, while row 20, 27
# and 28 stay untouched. Columns that are always blank for this block
# (I, Y, AA, AT, AY) are left alone.
# ---------------------------------------------------------------------------

# after-row -> source (before-row) mapping
$rowMap = @{
    17 = 23
    18 = 26
    19 = 21
    21 = 24
    22 = 17
    23 = 25
    24 = 29
    25 = 30
    26 = 22
    29 = 19
    30 = 18
}

# Columns present (non-blank) in every one of these rows.
$coreCols = @("A","B","C","D","E","F","G","H","P","Q","R","S","T","U","V","W","Z","AB","AD","AE","AG","AW","AX")

# Columns that are only present on some rows - copied when the source row
# has a value there, otherwise the destination cell is cleared (if it had
# something previously).
$optionalCols = @("M","AC","AI","AJ","AK","AL","AO")

$allCols = $coreCols + $optionalCols

function ColToNum($col) {
    $n = 0
    foreach ($ch in $col.ToCharArray()) {
        $n = $n * 26 + ([int][char]$ch - [int][char]'A' + 1)
    }
    return $n
}

# Snapshot the current (pre-edit) contents of every source row/column we
# need, BEFORE any writes happen (the mapping is a permutation / contains
# cycles, so we must not read-after-write).
$snapshot = @{}
foreach ($srcRow in $rowMap.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowData = @{}
        foreach ($col in $allCols) {
            $colNum = ColToNum $col
            $cell = $ws.Cells.Item($srcRow, $colNum)
            $rowData[$col] = $cell.Value2
        }
        $snapshot[$srcRow] = $rowData
    }
}

# Now write the snapshotted values into their destination rows.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $rowData = $snapshot[$srcRow]

    foreach ($col in $coreCols) {
        $colNum = ColToNum $col
        $ws.Cells.Item($destRow, $colNum).Value2 = $rowData[$col]
    }

    foreach ($col in $optionalCols) {
        $colNum = ColToNum $col
        $val = $rowData[$col]
        if ($null -ne $val -and $val -ne "") {
            $ws.Cells.Item($destRow, $colNum).Value2 = $val
        } else {
            $ws.Cells.Item($destRow, $colNum).ClearContents()
        }
    }
}
